$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update team-specific time-matrix probabilities (recomputed with additional data).
# Row 2
$ws.Range("B2").Value = 0.1956521739130435
$ws.Range("C2").Value = 0.5362318840579711
$ws.Range("J2").Value = 0.003623188405797101
$ws.Range("P2").Value = 0.177536231884058
$ws.Range("S2").Value = 0.08695652173913043

# Row 3
$ws.Range("B3").Value = 0.0130718954248366
$ws.Range("C3").Value = 0.0261437908496732
$ws.Range("J3").Value = 0.0196078431372549
$ws.Range("P3").Value = 0.7189542483660131

# Row 4
$ws.Range("J4").Value = 0.01923076923076923
$ws.Range("O4").Value = 0.01923076923076923
$ws.Range("P4").Value = 0.6730769230769231
$ws.Range("S4").Value = 0.2884615384615384

# Row 6
$ws.Range("B6").Value = 0.05860805860805861
$ws.Range("D6").Value = 0.01098901098901099
$ws.Range("F6").Value = 0.1025641025641026
$ws.Range("J6").Value = 0.1538461538461539
$ws.Range("O6").Value = 0.01831501831501832
$ws.Range("Q6").Value = 0.1684981684981685
$ws.Range("R6").Value = 0.1135531135531136
$ws.Range("S6").Value = 0.3736263736263736

# Row 7
$ws.Range("B7").Value = 0.132183908045977
$ws.Range("D7").Value = 0.02873563218390805
$ws.Range("F7").Value = 0.07471264367816093
$ws.Range("J7").Value = 0.1206896551724138
$ws.Range("O7").Value = 0.01724137931034483
$ws.Range("Q7").Value = 0.1149425287356322
$ws.Range("R7").Value = 0.08620689655172414
$ws.Range("S7").Value = 0.4252873563218391

# Row 8
$ws.Range("B8").Value = 0.07847533632286996
$ws.Range("D8").Value = 0.0179372197309417
$ws.Range("F8").Value = 0.06950672645739911
$ws.Range("J8").Value = 0.1390134529147982
$ws.Range("O8").Value = 0.02242152466367713
$ws.Range("Q8").Value = 0.1614349775784753
$ws.Range("R8").Value = 0.1098654708520179
$ws.Range("S8").Value = 0.4013452914798206

# Row 9
$ws.Range("B9").Value = 0.07352941176470588
$ws.Range("D9").Value = 0.02205882352941177
$ws.Range("F9").Value = 0.1102941176470588
$ws.Range("J9").Value = 0.1102941176470588
$ws.Range("O9").Value = 0.03308823529411765
$ws.Range("Q9").Value = 0.1985294117647059
$ws.Range("R9").Value = 0.1029411764705882
$ws.Range("S9").Value = 0.3492647058823529

# Row 10
$ws.Range("B10").Value = 0.0875643855776306
$ws.Range("D10").Value = 0.02281089036055924
$ws.Range("E10").Value = 0.0007358351729212656
$ws.Range("F10").Value = 0.09345106696100074
$ws.Range("J10").Value = 0.1074319352465048
$ws.Range("O10").Value = 0.01618837380426785
$ws.Range("Q10").Value = 0.1913171449595291
$ws.Range("R10").Value = 0.1103752759381898
$ws.Range("S10").Value = 0.3701250919793966

# Row 11
$ws.Range("J11").Value = 0.1118881118881119
$ws.Range("K11").Value = 0.1923076923076923
$ws.Range("L11").Value = 0.541958041958042
$ws.Range("S11").Value = 0.01748251748251748

# Row 12
$ws.Range("G12").Value = 0.6787878787878788
$ws.Range("J12").Value = 0.2424242424242424
$ws.Range("K12").Value = 0.006060606060606061
$ws.Range("L12").Value = 0.06060606060606061
$ws.Range("S12").Value = 0.01212121212121212

# Row 13
$ws.Range("F13").Value = 0.02173913043478261
$ws.Range("G13").Value = 0.6521739130434783
$ws.Range("J13").Value = 0.2826086956521739
$ws.Range("S13").Value = 0.04347826086956522

# Row 14
$ws.Range("G14").Value = 0.6
$ws.Range("S14").Value = 0.4

# Row 15
$ws.Range("F15").Value = 0.01700680272108844
$ws.Range("H15").Value = 0.1122448979591837
$ws.Range("I15").Value = 0.09863945578231292
$ws.Range("J15").Value = 0.391156462585034
$ws.Range("K15").Value = 0.05442176870748299
$ws.Range("M15").Value = 0.006802721088435374
$ws.Range("O15").Value = 0.09863945578231292
$ws.Range("S15").Value = 0.2210884353741497

# Row 16
$ws.Range("F16").Value = 0.01058201058201058
$ws.Range("H16").Value = 0.1534391534391534
$ws.Range("I16").Value = 0.1322751322751323
$ws.Range("J16").Value = 0.3703703703703703
$ws.Range("K16").Value = 0.1111111111111111
$ws.Range("M16").Value = 0.02645502645502645
$ws.Range("N16").Value = 0.005291005291005291
$ws.Range("O16").Value = 0.0582010582010582
$ws.Range("S16").Value = 0.1322751322751323

# Row 17
$ws.Range("F17").Value = 0.0111358574610245
$ws.Range("H17").Value = 0.1492204899777283
$ws.Range("I17").Value = 0.1158129175946548
$ws.Range("J17").Value = 0.4409799554565701
$ws.Range("K17").Value = 0.09799554565701558
$ws.Range("M17").Value = 0.022271714922049
$ws.Range("N17").Value = 0.004454342984409799
$ws.Range("O17").Value = 0.06013363028953229
$ws.Range("S17").Value = 0.09799554565701558

# Row 18
$ws.Range("F18").Value = 0.02592592592592593
$ws.Range("H18").Value = 0.1777777777777778
$ws.Range("I18").Value = 0.09259259259259259
$ws.Range("J18").Value = 0.4703703703703704
$ws.Range("K18").Value = 0.06296296296296296
$ws.Range("M18").Value = 0.01481481481481482
$ws.Range("N18").Value = 0.003703703703703704
$ws.Range("O18").Value = 0.07407407407407407
$ws.Range("S18").Value = 0.07777777777777778

# Row 19
$ws.Range("F19").Value = 0.006939090208172706
$ws.Range("H19").Value = 0.2074016962220509
$ws.Range("I19").Value = 0.1087124132613724
$ws.Range("J19").Value = 0.369313801079414
$ws.Range("K19").Value = 0.1010023130300694
$ws.Range("M19").Value = 0.02004626060138782
$ws.Range("N19").Value = 0.002313030069390902
$ws.Range("O19").Value = 0.09097918272937548
$ws.Range("S19").Value = 0.09329221279876639

Write-Output "Updated 113 cells in the team-specific matrix."